$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1: "Total Count:" -> "totalCount"
$ws.Range("C1").Formula = "totalCount"

# Row 2 data
$ws.Range("A2").Formula = "2025-02-23T13:05"
$ws.Range("B2").Formula = "SHIFT_1"

# Numeric-looking values must stay text (matching the original inline-string
# cells), so a leading apostrophe is used to force text entry instead of
# letting Excel coerce them into numbers.
$ws.Range("C2").Formula = "'2000"
$ws.Range("D2").Formula = "'0"
$ws.Range("E2").Formula = "'0"
$ws.Range("F2").Formula = "'0"
$ws.Range("G2").Formula = "'0"
$ws.Range("H2").Formula = "'0"
$ws.Range("I2").Formula = "'0"
$ws.Range("J2").Formula = "'0"
$ws.Range("K2").Formula = "'0"
$ws.Range("L2").Formula = "'0"
$ws.Range("M2").Formula = "'0"
$ws.Range("N2").Formula = "'0"
$ws.Range("O2").Formula = "'0"
$ws.Range("P2").Formula = "'0"
$ws.Range("Q2").Formula = "'0"
$ws.Range("R2").Formula = "'0"
$ws.Range("S2").Formula = "'0"
$ws.Range("T2").Formula = "'0"
$ws.Range("U2").Formula = "'0"
$ws.Range("V2").Formula = "'123"
$ws.Range("W2").Formula = "'321"
$ws.Range("X2").Formula = "'125"
